# Minor correction in files
$wb = $excel.ActiveWorkbook

# --- Cabling sheet: fix text in L2 and move selection there ---
$wsCabling = $wb.Worksheets.Item("Cabling")
$wsCabling.Range("L2").Value = "Cabling 3U4U3I3I3I"
$wsCabling.Range("L2").Select()

# --- DSPChannelMap sheet: update table values ---
$wsChannel = $wb.Worksheets.Item("DSPChannelMap")
$wsChannel.Range("D8").Value = 12
$wsChannel.Range("D9").Value = 14
$wsChannel.Range("B10").Value = 18
$wsChannel.Range("D10").Value = 15
$wsChannel.Range("B11").Value = 0
$wsChannel.Range("D11").Value = 16
$wsChannel.Range("B12").Value = 0
$wsChannel.Range("D12").Value = 17

# --- Make DSPChannelMap the active sheet/tab (was DSPFeederMap) ---
$wsChannel.Activate()
